# "cleaning of data completed"
#
# Data-cleaning pass on the country list workbook:
#   - "Sheet1" was a scratch/helper sheet used to spot duplicate Region
#     values while preparing the "LGBTIQ+" country table (helper column
#     with an IF() formula flagging non-unique rows). It is no longer
#     needed now that the cleaning is finished, so remove it entirely.
#   - Rename the remaining sheet from "LGBTIQ+" to "list" now that it is
#     a generic, cleaned country/region reference list.

$wb = $excel.ActiveWorkbook

# Avoid any "this sheet contains data, delete anyway?" prompt.
$excel.DisplayAlerts = $false

# Drop the helper de-duplication sheet.
$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining sheet to its new, final name.
$wb.Worksheets.Item("LGBTIQ+").Name = "list"

$excel.DisplayAlerts = $true
